# CARAGA_NEWCON.xlsx — "move final programs to separate folder" edit
#
# 1. Freeze the header row (row 1).
# 2. Set explicit column widths for columns A..AC (1..29).
# 3. Re-wrap / re-center the header row (A1:AB1): vertical=center, wrapText=1
#    (horizontal was already center).
# 4. Highlight the "Status as of ..." column (AC) in yellow, for both the
#    header (AC1) and the data cell (AC2), and bump the header's date text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Freeze panes: split/freeze after row 1 -------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.166666666666666
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668
$ws.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws.Columns.Item(7).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(8).ColumnWidth = 13.166666666666666
$ws.Columns.Item(9).ColumnWidth = 26.166666666666668
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 6.166666666666667
$ws.Columns.Item(12).ColumnWidth = 25.166666666666668
$ws.Columns.Item(13).ColumnWidth = 19.166666666666668
$ws.Columns.Item(14).ColumnWidth = 16.166666666666668
$ws.Columns.Item(15).ColumnWidth = 10.166666666666666
$ws.Columns.Item(16).ColumnWidth = 25.166666666666668
$ws.Columns.Item(17).ColumnWidth = 25.166666666666668
$ws.Columns.Item(18).ColumnWidth = 26.166666666666668
$ws.Columns.Item(19).ColumnWidth = 11.166666666666666
$ws.Columns.Item(20).ColumnWidth = 12.166666666666666
$ws.Columns.Item(21).ColumnWidth = 30.166666666666668
$ws.Columns.Item(22).ColumnWidth = 26.166666666666668
$ws.Columns.Item(23).ColumnWidth = 12.166666666666666
$ws.Columns.Item(24).ColumnWidth = 32.166666666666664
$ws.Columns.Item(25).ColumnWidth = 30.166666666666668
$ws.Columns.Item(26).ColumnWidth = 19.166666666666668
$ws.Columns.Item(27).ColumnWidth = 43.166666666666664
$ws.Columns.Item(28).ColumnWidth = 13.166666666666666
$ws.Columns.Item(29).ColumnWidth = 27.166666666666668

# --- Header row (A1:AB1): center vertically + wrap text ---------------------
$hdr = $ws.Range("A1:AB1")
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# --- "Status as of ..." column header (AC1): yellow fill, same alignment ----
$ac1 = $ws.Range("AC1")
$ac1.HorizontalAlignment = -4108
$ac1.VerticalAlignment = -4108
$ac1.WrapText = $true
$ac1.Interior.Color = 65535
$ac1.Value = "Status as of July 11, 2025"

# --- "Status as of ..." data cell (AC2): yellow fill -------------------------
$ws.Range("AC2").Interior.Color = 65535

Write-Host "edit applied"
